# Update scripts with new TPM values (NATMI LR-pairs output: Nmb-Nmbr.xlsx)
# The "Target cluster" column (D) for rows 2-7 is corrected from "Neutrophils"
# to "MuSCs", and the dependent expression/specificity metrics are updated
# to the newly recomputed TPM-based values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 2.401057
$ws.Range("H2").Value = 7.203171
$ws.Range("I2").Value = 0.1471169379718001
$ws.Range("J2").Value = 0.1471169379718
$ws.Range("M2").Value = 0.05595866666666666
$ws.Range("N2").Value = 0.167876
$ws.Range("Q2").Value = 0.1343599483106667
$ws.Range("R2").Value = 1.209239534796
$ws.Range("S2").Value = 0.1471169379718001
$ws.Range("T2").Value = 0.1471169379718

# Row 3
$ws.Range("D3").Value = "MuSCs"
$ws.Range("I3").Value = 0.5360701826106148
$ws.Range("J3").Value = 0.5360701826106148
$ws.Range("M3").Value = 0.05595866666666666
$ws.Range("N3").Value = 0.167876
$ws.Range("Q3").Value = 0.4895857881453333
$ws.Range("R3").Value = 4.406272093308
$ws.Range("S3").Value = 0.5360701826106148
$ws.Range("T3").Value = 0.5360701826106148

# Row 4
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 1.618009
$ws.Range("H4").Value = 4.854027
$ws.Range("I4").Value = 0.09913822524447119
$ws.Range("J4").Value = 0.09913822524447118
$ws.Range("M4").Value = 0.05595866666666666
$ws.Range("N4").Value = 0.167876
$ws.Range("Q4").Value = 0.09054162629466667
$ws.Range("R4").Value = 0.814874636652
$ws.Range("S4").Value = 0.09913822524447119
$ws.Range("T4").Value = 0.09913822524447118

# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.9788783333333333
$ws.Range("H5").Value = 2.936635
$ws.Range("I5").Value = 0.05997757781133019
$ws.Range("J5").Value = 0.05997757781133017
$ws.Range("M5").Value = 0.05595866666666666
$ws.Range("N5").Value = 0.167876
$ws.Range("Q5").Value = 0.05477672636222222
$ws.Range("R5").Value = 0.49299053726
$ws.Range("S5").Value = 0.05997757781133019
$ws.Range("T5").Value = 0.05997757781133017

# Row 6
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 2.140127
$ws.Range("H6").Value = 6.420381
$ws.Range("I6").Value = 0.1311293030989163
$ws.Range("J6").Value = 0.1311293030989162
$ws.Range("M6").Value = 0.05595866666666666
$ws.Range("N6").Value = 0.167876
$ws.Range("Q6").Value = 0.1197586534173333
$ws.Range("R6").Value = 1.077827880756
$ws.Range("S6").Value = 0.1311293030989163
$ws.Range("T6").Value = 0.1311293030989162

# Row 7
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.4336056666666666
$ws.Range("H7").Value = 1.300817
$ws.Range("I7").Value = 0.02656777326286756
$ws.Range("J7").Value = 0.02656777326286756
$ws.Range("M7").Value = 0.05595866666666666
$ws.Range("N7").Value = 0.167876
$ws.Range("Q7").Value = 0.02426399496577777
$ws.Range("R7").Value = 0.218375954692
$ws.Range("S7").Value = 0.02656777326286756
$ws.Range("T7").Value = 0.02656777326286756
